$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ToDO")

# Row 16 - "Code conventions opstellen" task gets completed details:
#   Werkelijke tijd (C16) = "15 min" (new shared string)
#   Voltooid (D16) = 5/3/2014, formatted like the other date cells (D15)
#   Solved (G16) = "Solved"
#   APP (H16) = "APP"

# Copy D15's number format onto D16 first so the new date reuses the
# existing date style instead of Excel auto-creating a new custom format.
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("C16").Value = "15 min"
$ws.Range("D16").Value = (Get-Date -Year 2014 -Month 5 -Day 3 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("G16").Value = "Solved"
$ws.Range("H16").Value = "APP"

# Update the active selection to match the post-edit state (C17).
$ws.Range("C17").Select()
